$wb = $excel.ActiveWorkbook

# --- Style change: the bold "title" font (14pt) and the bold "header row" font
#     are being unified into a single bold, white font used by both the
#     dashboard title (row 1) and the column header row (row 2) on every sheet.
foreach ($sht in $wb.Worksheets) {
    $sht.Range("A1").Font.Bold = $true
    $sht.Range("A1").Font.Size = 11
    $sht.Range("A1").Font.Color = 16777215

    $headerRow = $sht.Rows.Item(2)
    $headerRow.Font.Bold = $true
    $headerRow.Font.Color = 16777215
}

# --- Data changes on the "Training Dashboard" sheet, row 3 ---
$ws = $wb.Worksheets.Item("Training Dashboard")

# PERIOD TO EXPIRE
$ws.Range("H3").Value = 113

# LAST UPDATE -- keep this a literal text date (not an auto-converted serial)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
